# Daily attendance processing - 2026-01-21 22:40:25
# Reorders the "Recorded By" (column G) list so that the
# "dnasr281@gmail.com" entry no longer appears first in the
# comma-separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$targetName = "dnasr281@gmail.com"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -ne $value -and $value -like "*,*") {
        $parts = $value -split ",\s*"

        if ($parts.Count -eq 2 -and $parts[0].Trim() -eq $targetName) {
            $newValue = $parts[1].Trim() + ", " + $parts[0].Trim()
            $cell.Value2 = $newValue
        }
    }
}
